# Insert a new price record as row 194 in the "Mandarina" price sheet,
# pushing the existing rows 194:209 down to 195:210.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(194).Insert()

$ws.Cells.Item(194, 1).Value  = 1
$ws.Cells.Item(194, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(194, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(194, 4).Value  = 45223
$ws.Cells.Item(194, 5).Value  = 15
$ws.Cells.Item(194, 6).Value  = "Fruta"
$ws.Cells.Item(194, 7).Value  = 100102
$ws.Cells.Item(194, 8).Value  = "Cítricos"
$ws.Cells.Item(194, 9).Value  = 100102004
$ws.Cells.Item(194, 10).Value = "Mandarina"
$ws.Cells.Item(194, 11).Value = "Murcott"
$ws.Cells.Item(194, 12).Value = "Segunda"
$ws.Cells.Item(194, 13).Value = 300
$ws.Cells.Item(194, 14).Value = 14000
$ws.Cells.Item(194, 15).Value = 15000
$ws.Cells.Item(194, 16).Value = 14500
$ws.Cells.Item(194, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(194, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(194, 19).Value = 725
$ws.Cells.Item(194, 20).Value = 20
